$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.665.28'
$ws.Range('E2').Value = '  +0.47%  '
$ws.Range('D3').Value = '2.204.85'
$ws.Range('E3').Value = '  -1.82%  '
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '230.15'
$ws.Range('E5').Value = '  -0.76%  '
$ws.Range('E6').Value = '  -3.09%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '60.08'
$ws.Range('E7').Value = '  -6.15%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.401'
$ws.Range('E9').Value = '  -1.99%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '57.45'
$ws.Range('E10').Value = '  -3.68%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0888'
$ws.Range('E11').Value = '  -1.57%  '
$ws.Range('E12').Value = '  -1.34%  '
$ws.Range('D13').Value = '2.531.74'
$ws.Range('E13').Value = '  -1.93%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '15.35'
$ws.Range('E14').Value = '  -4.83%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '22.22'
$ws.Range('E15').Value = '  -1.52%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.69'
$ws.Range('E16').Value = '  +0.42%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.793'
$ws.Range('E17').Value = '  -3.83%  '
$ws.Range('D18').Value = '2.227.83'
$ws.Range('E18').Value = '  -0.77%  '
$ws.Range('D19').Value = '41.552.94'
$ws.Range('E19').Value = '  +0.48%  '
$ws.Range('D20').Value = '0.0₃0900'
$ws.Range('E20').Value = '  -4.59%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '71.93'
$ws.Range('E21').Value = '  -2.30%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.06'
$ws.Range('E22').Value = '  -1.55%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '242.96'
$ws.Range('E23').Value = '  -3.68%  '
$ws.Range('E24').Value = '  -0.11%  '
$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.35'
$ws.Range('E25').Value = '  -1.93%  '
$ws.Range('B26').Value = 'Toncoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.30'
$ws.Range('E26').Value = '  -1.28%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.66'
$ws.Range('E27').Value = '  -2.63%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '169.53'
$ws.Range('E28').Value = '  -2.13%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.139'
$ws.Range('E29').Value = '  -5.46%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '19.74'
$ws.Range('E30').Value = '  -3.41%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.42'
$ws.Range('E31').Value = '  -1.23%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.53'
$ws.Range('E32').Value = '  -10.02%  '
$ws.Range('E33').Value = '  -3.08%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.99'
$ws.Range('E34').Value = '  -0.90%  '
$ws.Range('E35').Value = '  -2.80%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0647'
$ws.Range('E36').Value = '  +1.49%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.45'
$ws.Range('E37').Value = '  -6.18%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.37'
$ws.Range('E38').Value = '  -3.19%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.53'
$ws.Range('E39').Value = '  -8.38%  '
$ws.Range('B40').Value = 'TerraClassic'
$ws.Range('C40').Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.000241'
$ws.Range('E40').Value = '  -9.83%  '
$ws.Range('B41').Value = 'BinanceUSD'
$ws.Range('C41').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.999'
$ws.Range('E41').Value = '  -0.24%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0237'
$ws.Range('E42').Value = '  -2.26%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.56'
$ws.Range('E43').Value = '  -4.61%  '
$ws.Range('B44').Value = 'Cronos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0956'
$ws.Range('E44').Value = '  -1.35%  '
$ws.Range('B45').Value = 'TrustWalletToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.20'
$ws.Range('E45').Value = '  -2.76%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '97.30'
$ws.Range('E46').Value = '  -5.16%  '
$ws.Range('D47').Value = '1.467.33'
$ws.Range('E47').Value = '  -2.88%  '
$ws.Range('E48').Value = '  -11.97%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '16.40'
$ws.Range('E49').Value = '  -7.56%  '
$ws.Range('E50').Value = '  -1.80%  '
$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.24'
$ws.Range('E51').Value = '  +4.42%  '
